$d = $word.ActiveDocument

# Locate the paragraph that ends with the "shell code" sentence (the one
# whose trailing whitespace needs trimming and after which the new
# "Determining the Parameters..." section is inserted).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*till it reaches our shell code.*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$r = $target.Range

# Strip any trailing whitespace that precedes the paragraph mark, leaving
# the sentence ending in "shell code." with no trailing spaces.
$paraMarkPos = $r.End - 1
$probeLen = 10
$probeStart = [Math]::Max($r.Start, $paraMarkPos - $probeLen)
$probe = $d.Range($probeStart, $paraMarkPos).Text
$trimCount = 0
for ($j = $probe.Length - 1; $j -ge 0; $j--) {
    if ($probe.Substring($j, 1) -eq " ") {
        $trimCount++
    } else {
        break
    }
}
if ($trimCount -gt 0) {
    $trimRange = $d.Range($paraMarkPos - $trimCount, $paraMarkPos)
    $trimRange.Text = ""
}

# Recompute the paragraph-mark position after trimming, then inject the two
# new paragraphs (a bold heading followed by a body paragraph) as raw
# WordprocessingML right after this paragraph.
$r = $target.Range
$insertPoint = $d.Range($r.End - 1, $r.End - 1)
$newParasXml = '<w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/></w:rPr><w:t>Determining the Parameters used in the Malicious Input</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t>Firstly, to figure out the buffer length we just keep feeding the program more and more data or you could create an absurdly large pattern using “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t>pattern_create</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t>” in</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve">the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t>metasploit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve"> framework. I used a pattern of size 3000 and sent it to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t>nweb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve"> server through </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t>ncat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve">. This created a core dump on our recreation of the environment which told us that there was a segmentation fault with signal 11. Using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t>pattern_offset</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:bCs/></w:rPr><w:t xml:space="preserve">, we figured out </w:t></w:r></w:p>'
$insertPoint.InsertXML($newParasXml) | Out-Null
